# Updated Flask app to be memory-friendly and Render-ready
#
# Data fix on the "Voters" sheet: a handful of RelativeName (column C)
# values get a "ः: " prefix added in front of the name.
#
#   Row 14 (VoterName दिनेश राम)        : सहदेव राम      -> ः: सहदेव राम
#   Row 20 (VoterName उमेश कुमार राम)    : सहदेव राम      -> ः: सहदेव राम
#   Row 18 (VoterName उमेश खतवे)        : सोमी खतवे      -> ः: सोमी खतवे
#   Row 49 (VoterName संतोष मंडल)       : चन्द्रदेव मंडल -> ः: चन्द्रदेव मंडल

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Value = "ः: सहदेव राम"
$ws.Range("C20").Value = "ः: सहदेव राम"
$ws.Range("C18").Value = "ः: सोमी खतवे"
$ws.Range("C49").Value = "ः: चन्द्रदेव मंडल"
